$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing "Chart" configuration row (row 11) and insert it
# right after the "Table" row (row 12), pushing nothing else around since
# it is appended at the end. Using Rows.Insert with a copy on the clipboard
# carries over the cell styles (so H/I keep the same style as the source
# row), matching how the sibling "Chart" row is formatted.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(13).Insert(-4121)

# The new row is an identical chart block except its "type" value is "Pie"
# instead of "Bar".
$ws.Range("D13").Value2 = "Pie"

# Update the active selection / scroll position recorded in the sheet view
# (also clears the previous topLeftCell scroll position).
$ws.Range("F19").Select()
